$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Dlk1"
$ws.Cells.Item(2, 3).Value = "Notch1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1167573333333333
$ws.Cells.Item(2, 8).Value = 0.350272
$ws.Cells.Item(2, 9).Value = 0.0006433000764991399
$ws.Cells.Item(2, 10).Value = 0.0006433000764991399
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 47.57896333333333
$ws.Cells.Item(2, 14).Value = 142.73689
$ws.Cells.Item(2, 15).Value = 0.450188452948237
$ws.Cells.Item(2, 16).Value = 0.4501884529482371
$ws.Cells.Item(2, 17).Value = 5.555192881564444
$ws.Cells.Item(2, 18).Value = 49.99673593407999
$ws.Cells.Item(2, 19).Value = 0.0002896062662206303
$ws.Cells.Item(2, 20).Value = 0.0002896062662206303

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Dlk1"
$ws.Cells.Item(3, 3).Value = "Notch1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1167573333333333
$ws.Cells.Item(3, 8).Value = 0.350272
$ws.Cells.Item(3, 9).Value = 0.0006433000764991399
$ws.Cells.Item(3, 10).Value = 0.0006433000764991399
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.829723666666666
$ws.Cells.Item(3, 14).Value = 29.489171
$ws.Cells.Item(3, 15).Value = 0.09300808131111737
$ws.Cells.Item(3, 16).Value = 0.09300808131111739
$ws.Cells.Item(3, 17).Value = 1.147692322723555
$ws.Cells.Item(3, 18).Value = 10.329230904512
$ws.Cells.Item(3, 19).Value = 0.00005983210582248003
$ws.Cells.Item(3, 20).Value = 0.00005983210582248004

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Dlk1"
$ws.Cells.Item(4, 3).Value = "Notch1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1167573333333333
$ws.Cells.Item(4, 8).Value = 0.350272
$ws.Cells.Item(4, 9).Value = 0.0006433000764991399
$ws.Cells.Item(4, 10).Value = 0.0006433000764991399
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 10.933664
$ws.Cells.Item(4, 14).Value = 32.800992
$ws.Cells.Item(4, 15).Value = 0.1034534789405002
$ws.Cells.Item(4, 16).Value = 0.1034534789405003
$ws.Cells.Item(4, 17).Value = 1.276585452202667
$ws.Cells.Item(4, 18).Value = 11.489269069824
$ws.Cells.Item(4, 19).Value = 0.00006655163091652597
$ws.Cells.Item(4, 20).Value = 0.00006655163091652598

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Dlk1"
$ws.Cells.Item(5, 3).Value = "Notch1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.1167573333333333
$ws.Cells.Item(5, 8).Value = 0.350272
$ws.Cells.Item(5, 9).Value = 0.0006433000764991399
$ws.Cells.Item(5, 10).Value = 0.0006433000764991399
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 37.34441866666666
$ws.Cells.Item(5, 14).Value = 112.033256
$ws.Cells.Item(5, 15).Value = 0.3533499868001453
$ws.Cells.Item(5, 16).Value = 0.3533499868001453
$ws.Cells.Item(5, 17).Value = 4.360234738403554
$ws.Cells.Item(5, 18).Value = 39.242112645632
$ws.Cells.Item(5, 19).Value = 0.0002273100735395035
$ws.Cells.Item(5, 20).Value = 0.0002273100735395036

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Dlk1"
$ws.Cells.Item(6, 3).Value = "Notch1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 174.184255
$ws.Cells.Item(6, 8).Value = 522.552765
$ws.Cells.Item(6, 9).Value = 0.9597062674131449
$ws.Cells.Item(6, 10).Value = 0.959706267413145
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 47.57896333333333
$ws.Cells.Item(6, 14).Value = 142.73689
$ws.Cells.Item(6, 15).Value = 0.450188452948237
$ws.Cells.Item(6, 16).Value = 0.4501884529482371
$ws.Cells.Item(6, 17).Value = 8287.506281888982
$ws.Cells.Item(6, 18).Value = 74587.55653700084
$ws.Cells.Item(6, 19).Value = 0.4320486798114507
$ws.Cells.Item(6, 20).Value = 0.4320486798114508

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Dlk1"
$ws.Cells.Item(7, 3).Value = "Notch1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 174.184255
$ws.Cells.Item(7, 8).Value = 522.552765
$ws.Cells.Item(7, 9).Value = 0.9597062674131449
$ws.Cells.Item(7, 10).Value = 0.959706267413145
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 9.829723666666666
$ws.Cells.Item(7, 14).Value = 29.489171
$ws.Cells.Item(7, 15).Value = 0.09300808131111737
$ws.Cells.Item(7, 16).Value = 0.09300808131111739
$ws.Cells.Item(7, 17).Value = 1712.183093734202
$ws.Cells.Item(7, 18).Value = 15409.64784360782
$ws.Cells.Item(7, 19).Value = 0.08926043855435073
$ws.Cells.Item(7, 20).Value = 0.08926043855435076

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Dlk1"
$ws.Cells.Item(8, 3).Value = "Notch1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 174.184255
$ws.Cells.Item(8, 8).Value = 522.552765
$ws.Cells.Item(8, 9).Value = 0.9597062674131449
$ws.Cells.Item(8, 10).Value = 0.959706267413145
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 10.933664
$ws.Cells.Item(8, 14).Value = 32.800992
$ws.Cells.Item(8, 15).Value = 0.1034534789405002
$ws.Cells.Item(8, 16).Value = 0.1034534789405003
$ws.Cells.Item(8, 17).Value = 1904.47211826032
$ws.Cells.Item(8, 18).Value = 17140.24906434288
$ws.Cells.Item(8, 19).Value = 0.09928495212489188
$ws.Cells.Item(8, 20).Value = 0.09928495212489191

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Dlk1"
$ws.Cells.Item(9, 3).Value = "Notch1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 174.184255
$ws.Cells.Item(9, 8).Value = 522.552765
$ws.Cells.Item(9, 9).Value = 0.9597062674131449
$ws.Cells.Item(9, 10).Value = 0.959706267413145
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 37.34441866666666
$ws.Cells.Item(9, 14).Value = 112.033256
$ws.Cells.Item(9, 15).Value = 0.3533499868001453
$ws.Cells.Item(9, 16).Value = 0.3533499868001453
$ws.Cells.Item(9, 17).Value = 6504.809743861426
$ws.Cells.Item(9, 18).Value = 58543.28769475284
$ws.Cells.Item(9, 19).Value = 0.3391121969224515
$ws.Cells.Item(9, 20).Value = 0.3391121969224515

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Dlk1"
$ws.Cells.Item(10, 3).Value = "Notch1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 7.196452999999999
$ws.Cells.Item(10, 8).Value = 21.589359
$ws.Cells.Item(10, 9).Value = 0.03965043251035593
$ws.Cells.Item(10, 10).Value = 0.03965043251035594
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 47.57896333333333
$ws.Cells.Item(10, 14).Value = 142.73689
$ws.Cells.Item(10, 15).Value = 0.450188452948237
$ws.Cells.Item(10, 16).Value = 0.4501884529482371
$ws.Cells.Item(10, 17).Value = 342.3997734170566
$ws.Cells.Item(10, 18).Value = 3081.59796075351
$ws.Cells.Item(10, 19).Value = 0.01785016687056562
$ws.Cells.Item(10, 20).Value = 0.01785016687056562

$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Dlk1"
$ws.Cells.Item(11, 3).Value = "Notch1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 7.196452999999999
$ws.Cells.Item(11, 8).Value = 21.589359
$ws.Cells.Item(11, 9).Value = 0.03965043251035593
$ws.Cells.Item(11, 10).Value = 0.03965043251035594
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 9.829723666666666
$ws.Cells.Item(11, 14).Value = 29.489171
$ws.Cells.Item(11, 15).Value = 0.09300808131111737
$ws.Cells.Item(11, 16).Value = 0.09300808131111739
$ws.Cells.Item(11, 17).Value = 70.73914437015432
$ws.Cells.Item(11, 18).Value = 636.6522993313889
$ws.Cells.Item(11, 19).Value = 0.003687810650944156
$ws.Cells.Item(11, 20).Value = 0.003687810650944157

$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Dlk1"
$ws.Cells.Item(12, 3).Value = "Notch1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 7.196452999999999
$ws.Cells.Item(12, 8).Value = 21.589359
$ws.Cells.Item(12, 9).Value = 0.03965043251035593
$ws.Cells.Item(12, 10).Value = 0.03965043251035594
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 10.933664
$ws.Cells.Item(12, 14).Value = 32.800992
$ws.Cells.Item(12, 15).Value = 0.1034534789405002
$ws.Cells.Item(12, 16).Value = 0.1034534789405003
$ws.Cells.Item(12, 17).Value = 78.683599093792
$ws.Cells.Item(12, 18).Value = 708.1523918441279
$ws.Cells.Item(12, 19).Value = 0.004101975184691834
$ws.Cells.Item(12, 20).Value = 0.004101975184691835

$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Dlk1"
$ws.Cells.Item(13, 3).Value = "Notch1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 7.196452999999999
$ws.Cells.Item(13, 8).Value = 21.589359
$ws.Cells.Item(13, 9).Value = 0.03965043251035593
$ws.Cells.Item(13, 10).Value = 0.03965043251035594
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 37.34441866666666
$ws.Cells.Item(13, 14).Value = 112.033256
$ws.Cells.Item(13, 15).Value = 0.3533499868001453
$ws.Cells.Item(13, 16).Value = 0.3533499868001453
$ws.Cells.Item(13, 17).Value = 268.7473537469893
$ws.Cells.Item(13, 18).Value = 2418.726183722904
$ws.Cells.Item(13, 19).Value = 0.01401047980415432
$ws.Cells.Item(13, 20).Value = 0.01401047980415432

